$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Insert the new "bibliography" slide at position 16 (right before the
# existing "Adaptation" slide), using the same "Title and Content" layout
# that the neighbouring slides use.
# ---------------------------------------------------------------------------
$newSlide = $p.Slides.Add(16, 2)

$newSlide.Shapes.Item(1).TextFrame.TextRange.Text = "bibliography"

$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "User Modeling: Recent Work, Prospects and Hazards1 – alfred kobsa.`rThe User Modeling Shell System BGP-MS - Alfred Kobsa and Wolfgang Pohl`rUser Modeling in Adaptive Hypermedia Educational Systems - António Constantino Martins, Luíz Faria, Carlos Vaz de Carvalho, Eurico Carrapatoso`rUser Modeling in Adaptive Interfaces- Pat Langley`r"

# Paragraph 1: "alfred kobsa." gets sz=2400
$body.Characters(54, 13).Font.Size = 24

# Paragraph 2: "Alfred Kobsa and Wolfgang Pohl" gets sz=2400
$body.Characters(108, 30).Font.Size = 24

# Paragraph 3: names after "Systems - " get sz=2400, the trailing space after
# "Carvalho," is bold (no size override), "Eurico Carrapatoso" gets sz=2000
$body.Characters(198, 64).Font.Size = 24
$body.Characters(262, 1).Font.Bold = $true
$body.Characters(263, 18).Font.Size = 20

# Paragraph 4: "Pat Langley" gets sz=2400
$body.Characters(320, 11).Font.Size = 24

# ---------------------------------------------------------------------------
# The former slide17 ("Question") shifts down to slide 18: mark it hidden and
# give it an (empty) transition element, matching the source deck.
# ---------------------------------------------------------------------------
$question = $p.Slides.Item(18)
$question.SlideShowTransition.Hidden = $true
$question.SlideShowTransition.SoundEffect.Name = ""
